$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2:K2").ClearContents()
$ws.Range("I3:K3").ClearContents()
$ws.Range("G4:K4").ClearContents()
$ws.Range("I5:K5").ClearContents()
$ws.Range("G6:K6").ClearContents()
$ws.Range("I7:K7").ClearContents()
$ws.Range("G8:K8").ClearContents()
$ws.Range("I9:K9").ClearContents()
$ws.Range("G10:K10").ClearContents()
$ws.Range("I11:K11").ClearContents()
$ws.Range("G12:K12").ClearContents()
$ws.Range("I13:K13").ClearContents()
$ws.Range("G14:K14").ClearContents()
$ws.Range("I15:K15").ClearContents()
$ws.Range("G16:K16").ClearContents()
$ws.Range("I17:K17").ClearContents()
$ws.Range("K18:K18").ClearContents()
$ws.Range("I19:K19").ClearContents()
$ws.Range("K20:K20").ClearContents()
$ws.Range("I21:K21").ClearContents()
$ws.Range("K22:K22").ClearContents()
$ws.Range("J23:K23").ClearContents()
$ws.Range("I24:K24").ClearContents()
$ws.Range("K26:K26").ClearContents()
$ws.Range("J27:K27").ClearContents()
$ws.Range("I28:K28").ClearContents()
$ws.Range("K30:K30").ClearContents()
$ws.Range("J31:K31").ClearContents()
$ws.Range("I32:K32").ClearContents()
$ws.Range("K34:K34").ClearContents()
$ws.Range("J35:K35").ClearContents()
$ws.Range("I36:K36").ClearContents()
$ws.Range("K38:K38").ClearContents()
$ws.Range("J39:K39").ClearContents()
$ws.Range("I40:K40").ClearContents()
$ws.Range("K42:K42").ClearContents()
$ws.Range("J43:K43").ClearContents()
$ws.Range("I44:J44").ClearContents()
